$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "31.410.57"
Set-TextValue "E2" "  +3.65%  "
Set-TextValue "D3" "2.009.87"
Set-TextValue "E3" "  +7.70%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "0.7906"
Set-TextValue "E5" "  +67.29%  "
Set-TextValue "D6" "260.20"
Set-TextValue "E6" "  +6.88%  "
Set-TextValue "D7" "0.9992"
Set-TextValue "E7" "  -0.11%  "
Set-TextValue "D8" "0.3613"
Set-TextValue "E8" "  +25.73%  "
Set-TextValue "D9" "28.73"
Set-TextValue "E9" "  +33.34%  "
Set-TextValue "D10" "0.07077"
Set-TextValue "E10" "  +9.20%  "
Set-TextValue "D11" "0.8505"
Set-TextValue "E11" "  +18.43%  "
Set-TextValue "D12" "0.08100"
Set-TextValue "E12" "  +4.06%  "
Set-TextValue "D13" "101.62"
Set-TextValue "E13" "  +5.00%  "
Set-TextValue "D14" "1.999.26"
Set-TextValue "E14" "  +7.11%  "
Set-TextValue "E15" "  +9.97%  "
Set-TextValue "D16" "276.52"
Set-TextValue "E16" "  -1.41%  "
Set-TextValue "D17" "31.388.77"
Set-TextValue "E17" "  +3.63%  "
Set-TextValue "D18" "14.71"
Set-TextValue "E18" "  +13.16%  "
Set-TextValue "D19" "5.955"
Set-TextValue "E19" "  +13.72%  "
Set-TextValue "D20" "0.000007947"
Set-TextValue "E20" "  +6.30%  "
Set-TextValue "D21" "2.264.36"
Set-TextValue "E21" "  +7.53%  "
Set-TextValue "D22" "0.9990"
Set-TextValue "E22" "  -0.11%  "
Set-TextValue "D23" "1.000"
Set-TextValue "E23" "  +0.02%  "
Set-TextValue "D24" "7.213"
Set-TextValue "E24" "  +15.37%  "
Set-TextValue "D25" "10.20"
Set-TextValue "E25" "  +13.67%  "
Set-TextValue "D26" "0.1515"
Set-TextValue "E26" "  +57.28%  "
Set-TextValue "D27" "164.70"
Set-TextValue "E27" "  +1.62%  "
Set-TextValue "D28" "20.11"
Set-TextValue "E28" "  +7.64%  "
Set-TextValue "E29" "  +28.23%  "
Set-TextValue "D30" "1.632"
Set-TextValue "E30" "  +10.06%  "
Set-TextValue "D31" "4.646"
Set-TextValue "E31" "  +10.33%  "
Set-TextValue "D32" "1.361"
Set-TextValue "E32" "  +3.66%  "
Set-TextValue "D33" "4.412"
Set-TextValue "E33" "  +7.26%  "
Set-TextValue "D34" "0.05236"
Set-TextValue "E34" "  +9.57%  "
Set-TextValue "D35" "1.226"
Set-TextValue "E35" "  +9.79%  "
Set-TextValue "D36" "0.7704"
Set-TextValue "E36" "  +12.71%  "
Set-TextValue "D37" "2.811"
Set-TextValue "E37" "  +3.70%  "
Set-TextValue "D38" "0.02014"
Set-TextValue "E38" "  +6.56%  "
Set-TextValue "D39" "2.954"
Set-TextValue "E39" "  +3.84%  "
Set-TextValue "D40" "81.31"
Set-TextValue "E40" "  +8.05%  "
Set-TextValue "D41" "6.716"
Set-TextValue "E41" "  +8.05%  "
Set-TextValue "D42" "2.200"
Set-TextValue "E42" "  +14.00%  "
Set-TextValue "D43" "0.4753"
Set-TextValue "E43" "  +13.25%  "
Set-TextValue "D44" "0.8604"
Set-TextValue "E44" "  +4.25%  "
Set-TextValue "D45" "104.93"
Set-TextValue "E45" "  +4.56%  "
Set-TextValue "D46" "1.000"
Set-TextValue "E46" "  +0.10%  "
Set-TextValue "B47" "Aptos"
Set-TextValue "C47" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.678"
Set-TextValue "E47" "  +10.38%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.935"
Set-TextValue "E48" "  +3.59%  "
Set-TextValue "D49" "0.4381"
Set-TextValue "E49" "  +13.09%  "
Set-TextValue "D50" "37.05"
Set-TextValue "E50" "  +6.04%  "
Set-TextValue "B51" "Algorand"
Set-TextValue "C51" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D51" "0.1189"
Set-TextValue "E51" "  +14.79%  "
